$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the "Incentives Available" column (old column H). Excel shifts
# everything after it one column to the left and drops the now-unused
# shared strings ("Incentives Available", "8,9", "1,12,21", "13,14,22",
# "4,23", "16,17").
$ws.Columns.Item(8).Delete()

# The worksheet carried a stale <sortState> (left over from sorting the old
# N3:P26 scratch columns on P). After the column shift those columns are now
# M3:O26 / O3:O26 - refresh the sort state to point at the new location
# (the range is effectively empty so this does not reorder any rows).
$sortRange = $ws.Range("M3:O26")
$keyRange = $ws.Range("O3:O26")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 2
$ws.Sort.Apply()

# Re-point the frozen pane / selection the same way the author left the file.
$ws.Range("L1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 5
